$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'247.04"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'0.93%"
$ws.Range("E2").Style = "Normal"
$ws.Range("E3").Value = "'5.25%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.067"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'1.21%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.05605"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-0.14%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'6.476"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-1.50%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.8131"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'0.28%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.8434"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'0.71%"
$ws.Range("E8").Style = "Normal"
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D9").Value = "'0.1338"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'0.00%"
$ws.Range("E9").Style = "Normal"
$ws.Range("B10").Value = "MandalaExchangeToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D10").Value = "'0.07000"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'0.71%"
$ws.Range("E10").Style = "Normal"
$ws.Range("B11").Value = "BitrueCoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D11").Value = "'0.02856"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'0.60%"
$ws.Range("E11").Style = "Normal"
$ws.Range("B12").Value = "BitMartToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D12").Value = "'0.09402"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-0.05%"
$ws.Range("E12").Style = "Normal"
$ws.Range("B13").Value = "BitForexToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D13").Value = "'0.001515"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.03%"
$ws.Range("E13").Style = "Normal"
$ws.Range("B14").Value = "One"
$ws.Range("C14").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D14").Value = "'0.0005959"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-93.83%"
$ws.Range("E14").Style = "Normal"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "'0.006164"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'1.34%"
$ws.Range("E15").Style = "Normal"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "'3.606"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'3.05%"
$ws.Range("E16").Style = "Normal"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").Value = "'3.013"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'0.24%"
$ws.Range("E17").Style = "Normal"
$ws.Range("B18").Value = "BTSEToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D18").Value = "'2.055"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-1.71%"
$ws.Range("E18").Style = "Normal"
$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D19").Value = "'0.3126"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'-2.21%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.03199"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-2.45%"
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'-1.36%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'3.738"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-0.11%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04652"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-0.53%"
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'-1.46%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.001243"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'0.03%"
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'1.45%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.00009600"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Value = "'0.0001397"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'-27.96%"
$ws.Range("E28").Style = "Normal"
$ws.Range("D40").Value = "'0.03667"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'0.21%"
$ws.Range("E40").Style = "Normal"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "'0.006150"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-1.28%"
$ws.Range("E41").Style = "Normal"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1058"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-22.40%"
$ws.Range("E42").Style = "Normal"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "'0.002500"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-8.60%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.008261"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'2.29%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005397"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'2.10%"
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'0.00%"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'-38.89%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.002600"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'27.36%"
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'0.00%"
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'0.00%"
$ws.Range("E50").Style = "Normal"
